# OLX Monitor 2026-02-19 09:32
#
# Appends a fresh monitoring snapshot (8 listings across profiles "poqui",
# "pokojewlublinie" and "dawnypatron") to the bottom of the running log on
# the PODSUMOWANIE sheet: rows 83:90, stamped 2026-02-19 09:32:21.
#
# The new block re-uses the formatting of the previous batch (rows 75:82,
# stamped 2026-02-19 08:46:44) which has the same 8 listings in the same
# order, so a straight Copy gives every new cell the right style before the
# values are overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$srcRange = $ws.Range("A75:H82")
$dstRange = $ws.Range("A83:H90")
$srcRange.Copy($dstRange)

$timestamp = "2026-02-19 09:32:21"

$rows = @(
    @{ Row = 83; Profile = "poqui";           Title = "Kawalerka po remoncie z funkcjonalną antresolą - ul. Jana Sawy";              Price = 2499;  Date = "28.10.2025"; Days = 113; Url = "https://www.olx.pl/d/oferta/kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger.html";                               Slug = "kawalerka-po-remoncie-z-funkcjonalna-antresola-ul-jana-sawy-CID3-ID183ger" },
    @{ Row = 84; Profile = "poqui";           Title = "Świeżo wykończone mieszkanie z dużym balkonem - Ponikwoda";                    Price = 2299;  Date = "19.01.2026"; Days = 30;  Url = "https://www.olx.pl/d/oferta/swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR.html";                                   Slug = "swiezo-wykonczone-mieszkanie-z-duzym-balkonem-ponikwoda-CID3-ID1951OR" },
    @{ Row = 85; Profile = "poqui";           Title = "Mieszkanie z KLIMATYZACJĄ 5 minut od UMCS, UP, KUL - Długosza";                Price = 2049;  Date = "19.12.2025"; Days = 61;  Url = "https://www.olx.pl/d/oferta/mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc.html";                               Slug = "mieszkanie-z-klimatyzacja-5-minut-od-umcs-up-kul-dlugosza-CID3-ID18KAEc" },
    @{ Row = 86; Profile = "poqui";           Title = "Przytulny pokój blisko Politechniki – ul. Przytulna";                          Price = 599;   Date = "10.10.2025"; Days = 132; Url = "https://www.olx.pl/d/oferta/przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz.html";                                         Slug = "przytulny-pokoj-blisko-politechniki-ul-przytulna-CID3-ID17NeTz" },
    @{ Row = 87; Profile = "pokojewlublinie"; Title = "WOLNY OD ZARAZ! Pokój jedynka, ul. Romanowskiego 58";                          Price = 58640; Date = "11.08.2025"; Days = 191; Url = "https://www.olx.pl/d/oferta/wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm.html";                                       Slug = "wolny-od-zaraz-pokoj-jedynka-ul-romanowskiego-58-CID3-ID16ZeYm" },
    @{ Row = 88; Profile = "pokojewlublinie"; Title = "WOLNY OD ZARAZ! Super lokalizacja, blisko centrum, ul. Paganiniego 12";        Price = 12640; Date = "19.01.2026"; Days = 30;  Url = "https://www.olx.pl/d/oferta/wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc.html";                       Slug = "wolny-od-zaraz-super-lokalizacja-blisko-centrum-ul-paganiniego-12-CID3-ID195dLc" },
    @{ Row = 89; Profile = "dawnypatron";     Title = "Ładny pokój jednoosobowy. Wynajmę duży pokój w centrum. ul Niecała 4.";        Price = 730;   Date = "20.09.2024"; Days = 516; Url = "https://www.olx.pl/d/oferta/ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM.html";                       Slug = "ladny-pokoj-jednoosobowy-wynajme-duzy-pokoj-w-centrum-ul-niecala-4-CID3-ID122jPM" },
    @{ Row = 90; Profile = "dawnypatron";     Title = "Mam do wynajęcia pokój dla os. pracującej lub studenta. Narutowicza 14";       Price = 14690; Date = "05.12.2025"; Days = 75;  Url = "https://www.olx.pl/d/oferta/mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv.html";                     Slug = "mam-do-wynajecia-pokoj-dla-os-pracujacej-lub-studenta-narutowicza-14-CID3-ID18ySfv" }
)

foreach ($r in $rows) {
    $row = $r.Row

    # Column A/E hold dot-separated "D.M.YYYY"-shaped (or dash/colon)
    # date-looking text that must stay literal text, not become a parsed
    # date serial. Force text mode for the write, then restore the
    # original "General" formatting below via a formats-only paste.
    $cA = $ws.Cells.Item($row, 1)
    $cA.NumberFormat = "@"
    $cA.Value = $timestamp

    $ws.Cells.Item($row, 2).Value = $r.Profile
    $ws.Cells.Item($row, 3).Value = $r.Title
    $ws.Cells.Item($row, 4).Value = $r.Price

    $cE = $ws.Cells.Item($row, 5)
    $cE.NumberFormat = "@"
    $cE.Value = $r.Date

    $ws.Cells.Item($row, 6).Value = $r.Days
    $ws.Cells.Item($row, 7).Value = $r.Url
    $ws.Cells.Item($row, 8).Value = $r.Slug
}

# Re-apply the original (General/left-aligned) number formatting that the
# "@" trick above temporarily overrode, without touching the values just
# written. Each destination column is re-formatted from its own source
# column only, so the paste can't bleed into neighbouring columns.
$ws.Range("A75:A82").Copy()
$ws.Range("A83:A90").PasteSpecial(-4122)
$ws.Range("E75:E82").Copy()
$ws.Range("E83:E90").PasteSpecial(-4122)
